$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 714503.4
$ws.Range("I2").Value = 714503.4
$ws.Range("J2").Value = 0.0
$ws.Range("K2").Value = 714503.4
$ws.Range("L2").Value = 0.0
$ws.Range("M2").Value = -714390.4
$ws.Range("N2").Value = ""
$ws.Range("H12").Value = 762.8182
$ws.Range("I12").Value = 1079.6
$ws.Range("K12").Value = 1079.6
$ws.Range("M12").Value = -909.5999999999999
$ws.Range("H40").Value = 4080.8572
$ws.Range("J40").Value = 5000.0
$ws.Range("L40").Value = 5000.0
$ws.Range("N40").Value = -5350.0
$ws.Range("H69").Value = 250013220.0
$ws.Range("I69").Value = 0.0
$ws.Range("K69").Value = 0.0
$ws.Range("M69").Value = ""
$ws.Range("H72").Value = 250013220.0
$ws.Range("I72").Value = 0.0
$ws.Range("K72").Value = 0.0
$ws.Range("M72").Value = ""
$ws.Range("H76").Value = 6666.6665
$ws.Range("I76").Value = 5000.0
$ws.Range("K76").Value = 5000.0
$ws.Range("M76").Value = -4685.0
$ws.Range("H79").Value = 6666.6665
$ws.Range("I79").Value = 5000.0
$ws.Range("K79").Value = 5000.0
$ws.Range("M79").Value = -3908.0
$ws.Range("H113").Value = 3482.389
$ws.Range("I113").Value = 2195.7144
$ws.Range("J113").Value = 4301.1816
$ws.Range("K113").Value = 2195.7144
$ws.Range("L113").Value = 4301.1816
$ws.Range("M113").Value = 1058.2856
$ws.Range("N113").Value = -10809.1816
$ws.Range("H116").Value = 5622.2856
$ws.Range("I116").Value = 6048.625
$ws.Range("J116").Value = 5053.8335
$ws.Range("K116").Value = 6048.625
$ws.Range("L116").Value = 5053.8335
$ws.Range("M116").Value = -2606.625
$ws.Range("N116").Value = -11937.8335
$ws.Range("H133").Value = 80000.0
$ws.Range("J133").Value = 80000.0
$ws.Range("L133").Value = 80000.0
$ws.Range("N133").Value = -90120.0
$ws.Range("H135").Value = 835.34784
$ws.Range("I135").Value = 835.34784
$ws.Range("K135").Value = 7518.130560000001
$ws.Range("M135").Value = -4983.130560000001
$ws.Range("H139").Value = 72998.4
$ws.Range("J139").Value = 72998.4
$ws.Range("L139").Value = 72998.4
$ws.Range("N139").Value = -83278.4

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 52500.0
$ws.Range("J76").Value = 52500.0
$ws.Range("L76").Value = 52500.0
$ws.Range("N76").Value = -53176.0
$ws.Range("H79").Value = 52500.0
$ws.Range("J79").Value = 52500.0
$ws.Range("L79").Value = 52500.0
$ws.Range("N79").Value = -54840.0

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 200000.0
$ws.Range("J70").Value = 200000.0
$ws.Range("L70").Value = 200000.0
$ws.Range("N70").Value = -200586.0
$ws.Range("H73").Value = 200000.0
$ws.Range("J73").Value = 200000.0
$ws.Range("L73").Value = 200000.0
$ws.Range("N73").Value = -202028.0
$ws.Range("H134").Value = 3302.1304
$ws.Range("I134").Value = 2049.9473
$ws.Range("K134").Value = 6149.841899999999
$ws.Range("M134").Value = -3614.841899999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 136.42857
$ws.Range("I35").Value = 136.42857
$ws.Range("K35").Value = 136.42857
$ws.Range("M35").Value = 157.57143
$ws.Range("H74").Value = 46138.0
$ws.Range("J74").Value = 46138.0
$ws.Range("L74").Value = 46138.0
$ws.Range("N74").Value = -47886.0
$ws.Range("H77").Value = 46138.0
$ws.Range("J77").Value = 46138.0
$ws.Range("L77").Value = 138414.0
$ws.Range("N77").Value = -147150.0
$ws.Range("H92").Value = 74973.0
$ws.Range("J92").Value = 74973.0
$ws.Range("L92").Value = 74973.0
$ws.Range("N92").Value = -79965.0
$ws.Range("H107").Value = 377.1111
$ws.Range("I107").Value = 428.14285
$ws.Range("J107").Value = 198.5
$ws.Range("K107").Value = 428.14285
$ws.Range("L107").Value = 198.5
$ws.Range("M107").Value = 1491.85715
$ws.Range("N107").Value = -4038.5
$ws.Range("H134").Value = 1899.0286
$ws.Range("I134").Value = 1898.8667
$ws.Range("K134").Value = 5696.6001
$ws.Range("M134").Value = -3161.6001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61.0
$ws.Range("J2").Value = 100.0
$ws.Range("L2").Value = 600.0
$ws.Range("N2").Value = -826.0
$ws.Range("H5").Value = 778.6429
$ws.Range("H38").Value = 1780.1666
$ws.Range("J38").Value = 2124.0
$ws.Range("L38").Value = 6372.0
$ws.Range("N38").Value = -7066.0
$ws.Range("H49").Value = 842.4
$ws.Range("I49").Value = 842.4
$ws.Range("K49").Value = 2527.2
$ws.Range("M49").Value = -2371.2
$ws.Range("H61").Value = 3362.6667
$ws.Range("I61").Value = 3362.6667
$ws.Range("K61").Value = 10088.0001
$ws.Range("M61").Value = -9873.000100000001
$ws.Range("H135").Value = 778.6429
$ws.Range("H138").Value = 2846.6
$ws.Range("I138").Value = 2725.0
$ws.Range("J138").Value = 3333.0
$ws.Range("K138").Value = 8175.0
$ws.Range("L138").Value = 9999.0
$ws.Range("M138").Value = -3035.0
$ws.Range("N138").Value = -20279.0
$ws.Range("H139").Value = 4568.0
$ws.Range("I139").Value = 4863.8184
$ws.Range("K139").Value = 14591.4552
$ws.Range("M139").Value = -9451.4552

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 7575.0
$ws.Range("J17").Value = 7575.0
$ws.Range("L17").Value = 7575.0
$ws.Range("N17").Value = -7911.0
$ws.Range("H44").Value = 5000.0
$ws.Range("J44").Value = 0.0
$ws.Range("L44").Value = 0.0
$ws.Range("N44").Value = ""
$ws.Range("H46").Value = 0.0
$ws.Range("J46").Value = 0.0
$ws.Range("L46").Value = 0.0
$ws.Range("N46").Value = ""
$ws.Range("H62").Value = 47000.0
$ws.Range("I62").Value = 44000.0
$ws.Range("K62").Value = 44000.0
$ws.Range("M62").Value = -43314.0
$ws.Range("H65").Value = 47000.0
$ws.Range("I65").Value = 44000.0
$ws.Range("K65").Value = 132000.0
$ws.Range("M65").Value = -128568.0
$ws.Range("H113").Value = 5574.143
$ws.Range("I113").Value = 3207.625
$ws.Range("K113").Value = 3207.625
$ws.Range("M113").Value = -1037.625
$ws.Range("H122").Value = 4223.6665
$ws.Range("I122").Value = 3020.4443
$ws.Range("K122").Value = 9061.332900000001
$ws.Range("M122").Value = -6611.332900000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 409.8
$ws.Range("J16").Value = 500.0
$ws.Range("L16").Value = 500.0
$ws.Range("N16").Value = -840.0
$ws.Range("H122").Value = 4289.0835
$ws.Range("I122").Value = 4605.926
$ws.Range("J122").Value = 3338.5557
$ws.Range("K122").Value = 13817.778
$ws.Range("L122").Value = 10015.6671
$ws.Range("M122").Value = -11367.778
$ws.Range("N122").Value = -14915.6671
$ws.Range("H132").Value = 1609.8334
$ws.Range("I132").Value = 1301.1818
$ws.Range("J132").Value = 5005.0
$ws.Range("K132").Value = 3903.5454
$ws.Range("L132").Value = 15015.0
$ws.Range("M132").Value = -1373.5454
$ws.Range("N132").Value = -20075.0

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1729038.0
$ws.Range("I2").Value = 2589187.8
$ws.Range("J2").Value = 8738.5
$ws.Range("K2").Value = 2589187.8
$ws.Range("L2").Value = 8738.5
$ws.Range("M2").Value = -2589075.8
$ws.Range("N2").Value = -8962.5
$ws.Range("H4").Value = 27750.25
$ws.Range("J4").Value = 27750.25
$ws.Range("L4").Value = 27750.25
$ws.Range("N4").Value = -27976.25
$ws.Range("H62").Value = 13827.821
$ws.Range("I62").Value = 5886.1665
$ws.Range("J62").Value = 15993.728
$ws.Range("K62").Value = 5886.1665
$ws.Range("L62").Value = 15993.728
$ws.Range("M62").Value = -5262.1665
$ws.Range("N62").Value = -17241.728
$ws.Range("H65").Value = 13827.821
$ws.Range("I65").Value = 5886.1665
$ws.Range("J65").Value = 15993.728
$ws.Range("K65").Value = 29430.8325
$ws.Range("L65").Value = 79968.64
$ws.Range("M65").Value = -26310.8325
$ws.Range("N65").Value = -86208.64
$ws.Range("H81").Value = 3568.6
$ws.Range("I81").Value = 2865.5557
$ws.Range("J81").Value = 4623.1665
$ws.Range("K81").Value = 5731.1114
$ws.Range("L81").Value = 9246.333
$ws.Range("M81").Value = -4670.1114
$ws.Range("N81").Value = -11368.333
$ws.Range("H84").Value = 3568.6
$ws.Range("I84").Value = 2865.5557
$ws.Range("J84").Value = 4623.1665
$ws.Range("K84").Value = 28655.557
$ws.Range("L84").Value = 46231.665
$ws.Range("M84").Value = -23351.557
$ws.Range("N84").Value = -56839.665
$ws.Range("H132").Value = 2517.4583
$ws.Range("I132").Value = 2242.5625
$ws.Range("K132").Value = 6727.6875
$ws.Range("M132").Value = -4197.6875
